# Apply updated matrix values for team_specific_matrix/North Dakota_B.xlsx
# (commit: "changes to team matrices from games pulled march 7")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1749502982107356
$ws.Range("C2").Value = 0.6003976143141153
$ws.Range("J2").Value = 0.009940357852882704
$ws.Range("P2").Value = 0.1351888667992048
$ws.Range("S2").Value = 0.07952286282306163
$ws.Range("B3").Value = 0.009708737864077669
$ws.Range("C3").Value = 0.02588996763754045
$ws.Range("J3").Value = 0.01941747572815534
$ws.Range("P3").Value = 0.7896440129449838
$ws.Range("S3").Value = 0.1553398058252427
$ws.Range("B6").Value = 0.05714285714285714
$ws.Range("D6").Value = 0.01666666666666667
$ws.Range("E6").Value = 0.002380952380952381
$ws.Range("F6").Value = 0.05714285714285714
$ws.Range("J6").Value = 0.2095238095238095
$ws.Range("O6").Value = 0.02380952380952381
$ws.Range("Q6").Value = 0.2023809523809524
$ws.Range("R6").Value = 0.05952380952380952
$ws.Range("S6").Value = 0.3714285714285714
$ws.Range("B7").Value = 0.08771929824561403
$ws.Range("D7").Value = 0.03258145363408521
$ws.Range("E7").Value = 0.002506265664160401
$ws.Range("F7").Value = 0.06265664160401002
$ws.Range("J7").Value = 0.08020050125313283
$ws.Range("O7").Value = 0.01754385964912281
$ws.Range("Q7").Value = 0.1929824561403509
$ws.Range("R7").Value = 0.09273182957393483
$ws.Range("S7").Value = 0.4310776942355889
$ws.Range("B8").Value = 0.08244422890397672
$ws.Range("D8").Value = 0.01939864209505335
$ws.Range("E8").Value = 0.0009699321047526673
$ws.Range("F8").Value = 0.05722599418040737
$ws.Range("J8").Value = 0.06692531522793405
$ws.Range("O8").Value = 0.01745877788554801
$ws.Range("Q8").Value = 0.2250242483026188
$ws.Range("R8").Value = 0.08632395732298739
$ws.Range("S8").Value = 0.4442289039767217
$ws.Range("B9").Value = 0.0949367088607595
$ws.Range("D9").Value = 0.0189873417721519
$ws.Range("E9").Value = 0.002109704641350211
$ws.Range("F9").Value = 0.05696202531645569
$ws.Range("J9").Value = 0.05696202531645569
$ws.Range("O9").Value = 0.02531645569620253
$ws.Range("Q9").Value = 0.1877637130801688
$ws.Range("R9").Value = 0.0759493670886076
$ws.Range("S9").Value = 0.4810126582278481
$ws.Range("B10").Value = 0.09440389294403893
$ws.Range("D10").Value = 0.02871046228710462
$ws.Range("E10").Value = 0.0009732360097323601
$ws.Range("F10").Value = 0.05693430656934306
$ws.Range("J10").Value = 0.09440389294403893
$ws.Range("O10").Value = 0.01800486618004866
$ws.Range("Q10").Value = 0.235036496350365
$ws.Range("R10").Value = 0.08467153284671533
$ws.Range("S10").Value = 0.3868613138686132
$ws.Range("G11").Value = 0.1360759493670886
$ws.Range("J11").Value = 0.07120253164556962
$ws.Range("K11").Value = 0.189873417721519
$ws.Range("L11").Value = 0.5474683544303798
$ws.Range("S11").Value = 0.05537974683544303
$ws.Range("G12").Value = 0.7486338797814208
$ws.Range("J12").Value = 0.1311475409836066
$ws.Range("K12").Value = 0.00819672131147541
$ws.Range("L12").Value = 0.0273224043715847
$ws.Range("S12").Value = 0.08469945355191257
$ws.Range("G13").Value = 0.5975609756097561
$ws.Range("J13").Value = 0.2804878048780488
$ws.Range("S13").Value = 0.1219512195121951
$ws.Range("F15").Value = 0.01659751037344398
$ws.Range("H15").Value = 0.1721991701244813
$ws.Range("I15").Value = 0.08091286307053942
$ws.Range("J15").Value = 0.2842323651452282
$ws.Range("K15").Value = 0.07261410788381743
$ws.Range("M15").Value = 0.006224066390041493
$ws.Range("N15").Value = 0.002074688796680498
$ws.Range("O15").Value = 0.06846473029045644
$ws.Range("S15").Value = 0.2966804979253112
$ws.Range("F16").Value = 0.007672634271099744
$ws.Range("H16").Value = 0.1918158567774936
$ws.Range("I16").Value = 0.1023017902813299
$ws.Range("J16").Value = 0.3529411764705883
$ws.Range("K16").Value = 0.1176470588235294
$ws.Range("M16").Value = 0.01023017902813299
$ws.Range("N16").Value = 0.002557544757033248
$ws.Range("O16").Value = 0.04603580562659847
$ws.Range("S16").Value = 0.1687979539641944
$ws.Range("F17").Value = 0.01043643263757116
$ws.Range("H17").Value = 0.198292220113852
$ws.Range("I17").Value = 0.1100569259962049
$ws.Range("J17").Value = 0.3206831119544592
$ws.Range("K17").Value = 0.09297912713472485
$ws.Range("M17").Value = 0.02087286527514232
$ws.Range("O17").Value = 0.07495256166982922
$ws.Range("S17").Value = 0.1717267552182163
$ws.Range("F18").Value = 0.0175
$ws.Range("H18").Value = 0.1525
$ws.Range("I18").Value = 0.12
$ws.Range("J18").Value = 0.3425
$ws.Range("K18").Value = 0.09
$ws.Range("M18").Value = 0.0225
$ws.Range("O18").Value = 0.08500000000000001
$ws.Range("S18").Value = 0.17
$ws.Range("F19").Value = 0.02691790040376851
$ws.Range("H19").Value = 0.2059219380888291
$ws.Range("I19").Value = 0.07873485868102288
$ws.Range("J19").Value = 0.2668236877523553
$ws.Range("K19").Value = 0.09320323014804845
$ws.Range("M19").Value = 0.01547779273216689
$ws.Range("N19").Value = 0.001009421265141319
$ws.Range("O19").Value = 0.06056527590847914
$ws.Range("S19").Value = 0.2513458950201884
